$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 4.083826541900635
$ws.Range("B1").Value = 4.010826110839844
$ws.Range("C1").Value = 2.616103172302246
$ws.Range("D1").Value = 2.221597671508789
$ws.Range("E1").Value = 1.767534494400024
